# Add traditional data generator comparison
# Updates the "Version 6" row (row 7) with new comparison values and
# moves the active selection, matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 ("Version 6") value updates ---
$ws.Cells.Item(7, 4).Value = 0.15    # D7: shift      0.1  -> 0.15
$ws.Cells.Item(7, 5).Value = $true   # E7: zoom       FALSE -> TRUE
$ws.Cells.Item(7, 6).Value = $false  # F7: flip       FALSE (unchanged)
$ws.Cells.Item(7, 7).Value = 98.7    # G7: ACC1       98.1 -> 98.7
$ws.Cells.Item(7, 8).Value = 98.1    # H7: ACC2       97.1 -> 98.1

# --- Update the active selection on the sheet ---
$ws.Range("F13").Select()

# --- Restore/update the workbook window layout ---
$win = $wb.Windows.Item(1)
$win.WindowState = -4143   # xlNormal
$win.Left = 1440
$win.Top = 1164
$win.Width = 21600
$win.Height = 11340
